# Auto-generated edit script: updates currentAveragePrice/profit columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    ,@("H33", 454.64)
    ,@("I33", 271.8)
    ,@("J33", 1186)
    ,@("K33", 271.8)
    ,@("L33", 1186)
    ,@("M33", -42.80000000000001)
    ,@("N33", -1644)
    ,@("H108", 77500)
    ,@("J108", 77500)
    ,@("L108", 77500)
    ,@("N108", -85180)
    ,@("H125", 2510.5)
    ,@("I125", 1649.75)
    ,@("J125", 3371.25)
    ,@("K125", 14847.75)
    ,@("L125", 30341.25)
    ,@("M125", -12387.75)
    ,@("N125", -35261.25)
    ,@("H137", 1781.9584)
    ,@("I137", 1059.96)
    ,@("K137", 3179.88)
    ,@("M137", -629.8800000000001)
    ,@("H138", 3610.5225)
    ,@("I138", 3729.7273)
    ,@("J138", 3552.2444)
    ,@("K138", 11189.1819)
    ,@("L138", 10656.7332)
    ,@("M138", -6049.1819)
    ,@("N138", -20936.7332)
)
foreach ($u in $updates) { $ws.Range($u[0]).Value = $u[1] }

$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    ,@("H2", 1117.0509)
    ,@("I2", 1086.5834)
    ,@("K2", 1086.5834)
    ,@("M2", -973.5834)
    ,@("H43", 90303)
    ,@("J43", 34995)
    ,@("L43", 34995)
    ,@("N43", -35621)
    ,@("H110", 1418.95)
    ,@("I110", 1255.3529)
    ,@("J110", 2346)
    ,@("K110", 1255.3529)
    ,@("L110", 2346)
    ,@("M110", 789.6470999999999)
    ,@("N110", -6436)
    ,@("H116", 1117.0509)
    ,@("I116", 1086.5834)
    ,@("K116", 1086.5834)
    ,@("M116", 1207.4166)
    ,@("H132", 601091.7)
    ,@("I132", 933485.0600000001)
    ,@("K132", 2800455.18)
    ,@("M132", -2797925.18)
)
foreach ($u in $updates) { $ws.Range($u[0]).Value = $u[1] }

$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    ,@("H3", 1117.0509)
    ,@("I3", 1086.5834)
    ,@("K3", 1086.5834)
    ,@("M3", -972.5834)
    ,@("H20", 2220.926)
    ,@("I20", 1956.7778)
    ,@("K20", 1956.7778)
    ,@("M20", -1709.7778)
    ,@("H105", 3226.1853)
    ,@("I105", 2221.111)
    ,@("J105", 5236.3335)
    ,@("K105", 2221.111)
    ,@("L105", 5236.3335)
    ,@("M105", -474.1109999999999)
    ,@("N105", -8730.333500000001)
    ,@("H107", 3264.2273)
    ,@("I107", 3357.762)
    ,@("J107", 1300)
    ,@("K107", 3357.762)
    ,@("L107", 1300)
    ,@("M107", -1437.762)
    ,@("N107", -5140)
    ,@("H134", 500453.62)
    ,@("I134", 748181)
    ,@("K134", 2244543)
    ,@("M134", -2242008)
)
foreach ($u in $updates) { $ws.Range($u[0]).Value = $u[1] }

$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    ,@("H31", 7858.3657)
    ,@("J31", 11415.8)
    ,@("L31", 11415.8)
    ,@("N31", -12005.8)
    ,@("H34", 7858.3657)
    ,@("J34", 11415.8)
    ,@("L34", 11415.8)
    ,@("N34", -11819.8)
    ,@("H122", 849.75)
    ,@("J122", 995)
    ,@("L122", 2985)
    ,@("N122", -7885)
    ,@("H132", 22761362)
    ,@("I132", 51897.715)
    ,@("K132", 155693.145)
    ,@("M132", -153163.145)
)
foreach ($u in $updates) { $ws.Range($u[0]).Value = $u[1] }

$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    ,@("H37", 105999.4)
    ,@("J37", 105999.4)
    ,@("L37", 317998.2)
    ,@("N37", -318222.2)
    ,@("H57", 8605.214)
    ,@("J57", 8614.5)
    ,@("L57", 25843.5)
    ,@("N57", -26961.5)
    ,@("H63", 5993.636)
    ,@("I63", 1995.8)
    ,@("J63", 9325.166999999999)
    ,@("K63", 5987.4)
    ,@("L63", 27975.501)
    ,@("M63", -5238.4)
    ,@("N63", -29473.501)
    ,@("H66", 5993.636)
    ,@("I66", 1995.8)
    ,@("J66", 9325.166999999999)
    ,@("K66", 17962.2)
    ,@("L66", 83926.503)
    ,@("M66", -14218.2)
    ,@("N66", -91414.503)
    ,@("H70", 6330)
    ,@("I70", 3000)
    ,@("J70", 7995)
    ,@("K70", 9000)
    ,@("L70", 23985)
    ,@("M70", -8685)
    ,@("N70", -24615)
    ,@("H73", 6330)
    ,@("I73", 3000)
    ,@("J73", 7995)
    ,@("K73", 9000)
    ,@("L73", 23985)
    ,@("M73", -7908)
    ,@("N73", -26169)
    ,@("H75", 4169.8)
    ,@("J75", 4042.1538)
    ,@("L75", 12126.4614)
    ,@("N75", -14122.4614)
    ,@("H76", 2699)
    ,@("I76", 2699)
    ,@("J76", 0)
    ,@("K76", 8097)
    ,@("L76", 0)
    ,@("M76", -7714)
    ,@("H78", 4169.8)
    ,@("J78", 4042.1538)
    ,@("L78", 36379.3842)
    ,@("N78", -46363.3842)
    ,@("H79", 2699)
    ,@("I79", 2699)
    ,@("J79", 0)
    ,@("K79", 8097)
    ,@("L79", 0)
    ,@("M79", -6771)
    ,@("H87", 20676.666)
    ,@("H90", 20676.666)
    ,@("H126", 7947.8)
    ,@("I126", 8022.25)
    ,@("K126", 24066.75)
    ,@("M126", -19126.75)
    ,@("H131", 10002.725)
    ,@("I131", 3472.375)
    ,@("J131", 12490.477)
    ,@("K131", 10417.125)
    ,@("L131", 37471.431)
    ,@("M131", -5377.125)
    ,@("N131", -47551.431)
)
foreach ($u in $updates) { $ws.Range($u[0]).Value = $u[1] }
$clears = @("N76", "N79")
foreach ($r in $clears) { $ws.Range($r).ClearContents() }

$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    ,@("H102", 2472.0908)
    ,@("I102", 1757.4)
    ,@("K102", 1757.4)
    ,@("M102", -135.4000000000001)
    ,@("H122", 3322.926)
    ,@("I122", 1814.174)
    ,@("J122", 11998.25)
    ,@("K122", 5442.522)
    ,@("L122", 35994.75)
    ,@("M122", -2992.522)
    ,@("N122", -40894.75)
)
foreach ($u in $updates) { $ws.Range($u[0]).Value = $u[1] }

$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    ,@("H7", 5969.8125)
    ,@("I7", 5701.1333)
    ,@("K7", 5701.1333)
    ,@("M7", -5589.1333)
    ,@("H40", 16020.6)
    ,@("J40", 2199.5)
    ,@("L40", 2199.5)
    ,@("N40", -2471.5)
    ,@("H104", 42722.25)
    ,@("J104", 42722.25)
    ,@("L104", 42722.25)
    ,@("N104", -49710.25)
    ,@("H126", 5969.8125)
    ,@("I126", 5701.1333)
    ,@("K126", 17103.3999)
    ,@("M126", -14633.3999)
    ,@("H136", 3959.7446)
    ,@("I136", 3407.3333)
    ,@("K136", 10221.9999)
    ,@("M136", -7671.999899999999)
)
foreach ($u in $updates) { $ws.Range($u[0]).Value = $u[1] }

$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    ,@("H37", 0)
    ,@("I37", 0)
    ,@("K37", 0)
    ,@("H122", 1789.3829)
    ,@("I122", 1322.9744)
    ,@("K122", 3968.9232)
    ,@("M122", -1518.9232)
    ,@("H136", 9095174)
    ,@("I136", 9873932)
    ,@("K136", 29621796)
    ,@("M136", -29619246)
)
foreach ($u in $updates) { $ws.Range($u[0]).Value = $u[1] }
$clears = @("M37")
foreach ($r in $clears) { $ws.Range($r).ClearContents() }

Write-Output "Updated $([string]203) cells across 8 sheets."